$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F2").Value = 2.68
$ws.Range("I2").Value = 2.7
$ws.Range("N2").Value = 5.2
$ws.Range("O2").Value = 1.21
$ws.Range("Q2").Value = 1.66
$ws.Range("R2").Value = 1.56
$ws.Range("S2").Value = 2.62
$ws.Range("T2").Value = 1.57
$ws.Range("V2").Value = 1.58
$ws.Range("L3").Value = 1.33
$ws.Range("I4").Value = 2.48
$ws.Range("K4").Value = 3.75
$ws.Range("V4").Value = 1.68
$ws.Range("W4").Value = 1.41
$ws.Range("Z4").Value = 19
$ws.Range("AJ4").Value = 70
$ws.Range("N5").Value = 4.2
$ws.Range("Q5").Value = 1.5
$ws.Range("R6").Value = 1.66
$ws.Range("S6").Value = 2.3
$ws.Range("F7").Value = 2.4
$ws.Range("I7").Value = 3.55
$ws.Range("N7").Value = 3.4
$ws.Range("Q7").Value = 1.99
$ws.Range("V7").Value = 1.4
$ws.Range("F8").Value = 3.2
$ws.Range("G8").Value = 4.1
$ws.Range("H8").Value = 2.02
$ws.Range("I8").Value = 2.26
$ws.Range("J8").Value = 3.5
$ws.Range("K8").Value = 5.2
$ws.Range("N8").Value = 4.6
$ws.Range("P8").Value = 2.22
$ws.Range("Q8").Value = 1.55
$ws.Range("R8").Value = 1.58
$ws.Range("S8").Value = 2.2
$ws.Range("U8").Value = 2.48
$ws.Range("V8").Value = 1.79
$ws.Range("I9").Value = 2.68
$ws.Range("P9").Value = 1.49
$ws.Range("V9").Value = 1.59
$ws.Range("AO9").Value = 21
$ws.Range("F10").Value = 2.8
$ws.Range("G10").Value = 2.94
$ws.Range("H10").Value = 2.68
$ws.Range("I10").Value = 2.84
$ws.Range("K10").Value = 3.6
$ws.Range("N10").Value = 3.6
$ws.Range("Q10").Value = 2
$ws.Range("R10").Value = 1.35
$ws.Range("S10").Value = 3.55
$ws.Range("W10").Value = 1.52
$ws.Range("Y10").Value = 11.5
$ws.Range("Z10").Value = 18
$ws.Range("AA10").Value = 42
$ws.Range("AB10").Value = 12
$ws.Range("AC10").Value = 8
$ws.Range("AD10").Value = 13
$ws.Range("AE10").Value = 32
$ws.Range("AF10").Value = 19
$ws.Range("AG10").Value = 13
$ws.Range("AH10").Value = 18
$ws.Range("AK10").Value = 32
$ws.Range("AM10").Value = 130
$ws.Range("AN10").Value = 27
$ws.Range("AO10").Value = 30
$ws.Range("F11").Value = 4.2
$ws.Range("I11").Value = 2.1
$ws.Range("K11").Value = 4.2
$ws.Range("T11").Value = 2.06
$ws.Range("V11").Value = 1.91
$ws.Range("W11").Value = 1.21
$ws.Range("AB11").Value = 15
$ws.Range("G12").Value = 3.9
$ws.Range("I12").Value = 2.56
$ws.Range("N12").Value = 2.16
$ws.Range("O12").Value = 1.74
$ws.Range("P12").Value = 1.39
$ws.Range("U12").Value = 1.58
$ws.Range("V12").Value = 1.64
$ws.Range("W12").Value = 1.35
$ws.Range("AO12").Value = 60
$ws.Range("G13").Value = 2.22
$ws.Range("H13").Value = 4.4
$ws.Range("I13").Value = 5.7
$ws.Range("J13").Value = 3.1
$ws.Range("N13").Value = 2.48
$ws.Range("V13").Value = 1.22
$ws.Range("W13").Value = 1.83
$ws.Range("K15").Value = 3.85
